# Fruta / hortaliza, semanal
# Insert a new week's worth of price rows (Primera/Segunda/Tercera) for
# "Agrícola del Norte S.A. de Arica" - Tomate, at the top of the data
# block (just below the header row), pushing the existing data down by
# three rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three new rows starting at row 911; everything from 911 down
# (previously ending at 933) shifts to 914..936.
$ws.Range("A911:R913").EntireRow.Insert()

# New weekly data block: Date 2022-07-05 (serial 44747).
$rows = @(
    @{ Row = 911; Cal = "Primera"; J = 300; K = 3000; L = 3500; M = 3250; P = 325 },
    @{ Row = 912; Cal = "Segunda"; J = 350; K = 2500; L = 3000; M = 2750; P = 275 },
    @{ Row = 913; Cal = "Tercera"; J = 400; K = 2000; L = 2500; M = 2250; P = 225 }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = 1
    $ws.Cells.Item($row, 2).Value = "Agrícola del Norte S.A. de Arica"
    $ws.Cells.Item($row, 3).Value = "Arica y Parinacota"
    $ws.Cells.Item($row, 4).Value = 44747
    $ws.Cells.Item($row, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($row, 5).Value = 15
    $ws.Cells.Item($row, 6).Value = 100112020
    $ws.Cells.Item($row, 7).Value = "Tomate"
    $ws.Cells.Item($row, 8).Value = "Larga vida"
    $ws.Cells.Item($row, 9).Value = $r.Cal
    $ws.Cells.Item($row, 10).Value = $r.J
    $ws.Cells.Item($row, 11).Value = $r.K
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
    $ws.Cells.Item($row, 14).Value = "`$/caja 10 kilos"
    $ws.Cells.Item($row, 15).Value = "Región de Arica y Parinacota"
    $ws.Cells.Item($row, 16).Value = $r.P
    $ws.Cells.Item($row, 17).Value = 10
    $ws.Cells.Item($row, 18).Value = "Hortaliza"
}
